$wb = $excel.ActiveWorkbook

# --- Update the conversion message on sheet "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.36 = 12818.79 pesos`n✅ 12818.79 pesos = 3.36 = 976.58 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("O10").Value = 3820
$ws2.Range("N12").Value = 3820
$ws2.Range("O12").Value = 291.02
